# Add Denmark, Sweden and Norway market test-data sheets, cloned from the
# existing "Belgium" sheet (same layout/styles), each with its own
# Market / NGC reference strings in B2 / B4.
# Also updates tab selection: UK is no longer the active tab, the
# newly-added "Norway" sheet becomes the active tab, and a few sheet-level
# selections are refreshed to match.

$wb = $excel.ActiveWorkbook

$belgium = $wb.Worksheets.Item("Belgium")

# --- Denmark -----------------------------------------------------------
$belgium.Copy($null, $belgium)
$denmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Name = "Denmark"
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").Value = "NGC-3446/T2009"

# --- Sweden --------------------------------------------------------------
$denmark.Copy($null, $denmark)
$sweden = $wb.Worksheets.Item($wb.Worksheets.Count)
$sweden.Name = "Sweden"
$sweden.Range("B2").Value = "Sweden Market"
$sweden.Range("B4").Value = "NGC-3465/T2021"

# --- Norway --------------------------------------------------------------
$sweden.Copy($null, $sweden)
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"
$norway.Range("B2").Value = "Norway Market"
$norway.Range("B4").Value = "NGC-3464/T1924"

# --- Selections ----------------------------------------------------------
# UK is no longer the selected tab; its cursor moved to B19.
$uk = $wb.Worksheets.Item("UK")
$uk.Activate()
$uk.Range("B19").Select()

# Denmark / Sweden have their whole grid selected (not the active tab).
$denmark.Activate()
$denmark.Cells.Select()

$sweden.Activate()
$sweden.Cells.Select()

# Norway ends up the active tab, with B2:B4 selected.
$norway.Activate()
$norway.Range("B2:B4").Select()
